# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Especial / Primera, "$/caja 18 kilos
# granel") above the existing Membrillo - Feria Lagunitas de Puerto Montt
# records, pushing the old rows 65-73 down to 67-75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: insert 2 blank rows at 65, shifting rows 65:73 down to 67:75.
$ws.Rows("65:66").Insert()

# New row 65: Especial, $/caja 18 kilos granel
$ws.Range("A65").Value = 4
$ws.Range("B65").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C65").Value = "Los Lagos"
$ws.Range("D65").Value = 44642
$ws.Range("E65").Value = 10
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100104
$ws.Range("H65").Value = "Frutos de pepita"
$ws.Range("I65").Value = 100104003
$ws.Range("J65").Value = "Membrillo"
$ws.Range("K65").Value = "Champion"
$ws.Range("L65").Value = "Especial"
$ws.Range("M65").Value = 300
$ws.Range("N65").Value = 19000
$ws.Range("O65").Value = 19000
$ws.Range("P65").Value = 19000
$ws.Range("Q65").Value = '$/caja 18 kilos granel'
$ws.Range("R65").Value = "Región de O'Higgins"
$ws.Range("S65").Value = 1056
$ws.Range("T65").Value = 18

# New row 66: Primera, $/caja 18 kilos granel
$ws.Range("A66").Value = 4
$ws.Range("B66").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C66").Value = "Los Lagos"
$ws.Range("D66").Value = 44642
$ws.Range("E66").Value = 10
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100104
$ws.Range("H66").Value = "Frutos de pepita"
$ws.Range("I66").Value = 100104003
$ws.Range("J66").Value = "Membrillo"
$ws.Range("K66").Value = "Champion"
$ws.Range("L66").Value = "Primera"
$ws.Range("M66").Value = 300
$ws.Range("N66").Value = 15000
$ws.Range("O66").Value = 15000
$ws.Range("P66").Value = 15000
$ws.Range("Q66").Value = '$/caja 18 kilos granel'
$ws.Range("R66").Value = "Región de O'Higgins"
$ws.Range("S66").Value = 833
$ws.Range("T66").Value = 18
